# Generate Report for Handoff
# Update status + timestamps for the zh-cn / de-de localization rows now that
# the handoff package has been generated, and widen the status columns to
# fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
# Status columns for zh-cn (E) and de-de (F), and the "Latest HO Xliff
# Generate Date" column (G).
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-29 12:58:59"

# --- zh-cn sheet ------------------------------------------------------------
# Status column (C) and Latest Handoff Datetime column (H).
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-29 12:58:55"

# --- de-de sheet ------------------------------------------------------------
# Status column (C) and Latest Handoff Datetime column (H).
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-29 12:58:59"

# --- Widen status columns to fit the new "Ready for handoff" text ----------
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # F: de-de status
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33        # C: Status
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33        # C: Status
